$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns C:D (full/tipo/link shift to E/F/G)
$ws.Range("C:D").Insert()

# Headers
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Per-row data: modelo (C), politica (D), tipo (F, lowercase), link (G, updated position/tracking_id)
$ws.Range("C2").Value = "FONTE 200A LITE"
$ws.Range("D2").Value = "Igual"
$ws.Range("F2").Value = "premium"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:1056404169#searchVariation=MLB24154371&position=3&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C3").Value = "FONTE 200 BOB"
$ws.Range("D3").Value = "Igual"
$ws.Range("F3").Value = "classico"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:1056404169#searchVariation=MLB24834408&position=1&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C4").Value = "Sem Modelo"
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27685629?pdp_filters=seller_id:1056404169#searchVariation=MLB27685629&position=2&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C5").Value = "FONTE 90 BOB"
$ws.Range("D5").Value = "Igual"
$ws.Range("F5").Value = "classico"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:1056404169#searchVariation=MLB21562641&position=9&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C6").Value = "FONTE 70A"
$ws.Range("D6").Value = "Igual"
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:1056404169#searchVariation=MLB21455208&position=4&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C7").Value = "FONTE 120A LITE"
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-120a-storm-lite-12v-bivolt-cor-preto/p/MLB23998473?pdp_filters=seller_id:1056404169#searchVariation=MLB23998473&position=7&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C8").Value = "FONTE 40A"
$ws.Range("D8").Value = "Igual"
$ws.Range("F8").Value = "classico"
$ws.Range("G8").Value = "https://www.mercadolivre.com.br/fonte-automotiva-40-amperes-jfa-storm-red-line-cca-sci-smart-cor-preto/p/MLB21621306?pdp_filters=seller_id:1056404169#searchVariation=MLB21621306&position=10&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C9").Value = "FONTE 40A"
$ws.Range("D9").Value = "Igual"
$ws.Range("F9").Value = "premium"
$ws.Range("G9").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-storm-40a-bivolt-12v-cor-preto/p/MLB22569833?pdp_filters=seller_id:1056404169#searchVariation=MLB22569833&position=8&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C10").Value = "Sem Modelo"
$ws.Range("F10").Value = "classico"
$ws.Range("G10").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:1056404169#searchVariation=MLB27687422&position=6&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C11").Value = "FONTE 120 BOB"
$ws.Range("D11").Value = "Igual"
$ws.Range("F11").Value = "classico"
$ws.Range("G11").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto/p/MLB22144397?pdp_filters=seller_id:1056404169#searchVariation=MLB22144397&position=16&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C12").Value = "FONTE 60A LITE"
$ws.Range("D12").Value = "Igual"
$ws.Range("F12").Value = "classico"
$ws.Range("G12").Value = "https://www.mercadolivre.com.br/jfa-fonte-carregador-storm-lite-60a-3000-w-preto/p/MLB23456525?pdp_filters=seller_id:1056404169#searchVariation=MLB23456525&position=17&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C13").Value = "Sem Modelo"
$ws.Range("F13").Value = "classico"
$ws.Range("G13").Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-azul/p/MLB28722231?pdp_filters=seller_id:1056404169#searchVariation=MLB28722231&position=18&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C14").Value = "FONTE 70A"
$ws.Range("D14").Value = "Igual"
$ws.Range("F14").Value = "premium"
$ws.Range("G14").Value = "https://www.mercadolivre.com.br/fonte-jfa-storm-modelo-com-70-amperes-para-carro/p/MLB27622275?pdp_filters=seller_id:1056404169#searchVariation=MLB27622275&position=19&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C15").Value = "FONTE 200A"
$ws.Range("D15").Value = "Igual"
$ws.Range("F15").Value = "premium"
$ws.Range("G15").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotiva-jfa-200a-slim-bivolt-voltimetro/p/MLB21348561?pdp_filters=seller_id:1056404169#searchVariation=MLB21348561&position=15&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C16").Value = "Sem Modelo"
$ws.Range("F16").Value = "classico"
$ws.Range("G16").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-k1200-pretolaranja-1200mt/p/MLB28357019?pdp_filters=seller_id:1056404169#searchVariation=MLB28357019&position=20&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C17").Value = "FONTE 60A"
$ws.Range("D17").Value = "Igual"
$ws.Range("F17").Value = "classico"
$ws.Range("G17").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:1056404169#searchVariation=MLB21320712&position=14&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C18").Value = "Sem Modelo"
$ws.Range("F18").Value = "classico"
$ws.Range("G18").Value = "https://www.mercadolivre.com.br/kit-controle-longa-distncia-jfa-k1200-completo-preto/p/MLB29770584?pdp_filters=seller_id:1056404169#searchVariation=MLB29770584&position=12&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C19").Value = "FONTE 120A"
$ws.Range("D19").Value = "Igual"
$ws.Range("F19").Value = "classico"
$ws.Range("G19").Value = "https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:1056404169#searchVariation=MLB21392652&position=13&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C20").Value = "FONTE 90 BOB"
$ws.Range("D20").Value = "Igual"
$ws.Range("F20").Value = "premium"
$ws.Range("G20").Value = "https://produto.mercadolivre.com.br/MLB-3863782558-nova-fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-_JM#position%3D21%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C21").Value = "FONTE 120A LITE"
$ws.Range("D21").Value = "Igual"
$ws.Range("F21").Value = "premium"
$ws.Range("G21").Value = "https://produto.mercadolivre.com.br/MLB-4131026512-fonte-automotiva-jfa-storm-lite-120a-bivolt-carregador-som-_JM#position%3D22%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C22").Value = "Sem Modelo"
$ws.Range("F22").Value = "classico"
$ws.Range("G22").Value = "https://produto.mercadolivre.com.br/MLB-3789324528-controle-longa-distancia-jfa-redline-wr-p-aparelho-original-_JM#position%3D23%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C23").Value = "FONTE 70A LITE"
$ws.Range("D23").Value = "Igual"
$ws.Range("F23").Value = "premium"
$ws.Range("G23").Value = "https://produto.mercadolivre.com.br/MLB-3350295935-fonte-automotiva-jfa-storm-lite-70a-bivolt-carregador-_JM#position%3D24%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C24").Value = "FONTE 70A LITE"
$ws.Range("D24").Value = "Igual"
$ws.Range("F24").Value = "classico"
$ws.Range("G24").Value = "https://produto.mercadolivre.com.br/MLB-3707207616-fonte-automotiva-jfa-storm-lite-70a-bivolt-carregador-som-_JM#position%3D25%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C25").Value = "Sem Modelo"
$ws.Range("F25").Value = "premium"
$ws.Range("G25").Value = "https://produto.mercadolivre.com.br/MLB-3344087007-controle-longa-distancia-jfa-k1200-azul-completo-o-melhor-_JM#position%3D26%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C26").Value = "Sem Modelo"
$ws.Range("F26").Value = "premium"
$ws.Range("G26").Value = "https://produto.mercadolivre.com.br/MLB-3693932379-controle-jfa-k1200-preto-com-laranja-longa-distncia-1200mt-_JM#position%3D27%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C27").Value = "Sem Modelo"
$ws.Range("F27").Value = "classico"
$ws.Range("G27").Value = "https://produto.mercadolivre.com.br/MLB-3572083997-controle-jfa-k1200-preto-com-laranja-longa-distncia-1200mt-_JM#position%3D28%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C28").Value = "FONTE 200A LITE"
$ws.Range("D28").Value = "Igual"
$ws.Range("F28").Value = "premium"
$ws.Range("G28").Value = "https://produto.mercadolivre.com.br/MLB-3724493418-fonte-automotiva-jfa-storm-lite-200a-bivolt-carregador-_JM#position%3D29%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C29").Value = "Sem Modelo"
$ws.Range("F29").Value = "premium"
$ws.Range("G29").Value = "https://produto.mercadolivre.com.br/MLB-3497393217-controle-jfa-acqua-prova-dagua-longa-distncia-1200m-branco-_JM#position%3D30%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C30").Value = "FONTE 120A"
$ws.Range("D30").Value = "Igual"
$ws.Range("F30").Value = "premium"
$ws.Range("G30").Value = "https://produto.mercadolivre.com.br/MLB-3706869734-fonte-automotiva-carregador-bateria-jfa-storm-120a-amperes-_JM#position%3D31%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C31").Value = "FONTE 200 BOB"
$ws.Range("D31").Value = "Igual"
$ws.Range("F31").Value = "premium"
$ws.Range("G31").Value = "https://produto.mercadolivre.com.br/MLB-3344001937-fonte-carregador-jfa-bob-storm-200a-bivolt-_JM?searchVariation=178756761911#searchVariation%3D178756761911%26position%3D32%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C32").Value = "Sem Modelo"
$ws.Range("F32").Value = "premium"
$ws.Range("G32").Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-preto/p/MLB28687615?pdp_filters=seller_id:1056404169#searchVariation=MLB28687615&position=11&search_layout=stack&type=product&tracking_id=4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C33").Value = "FONTE 60A"
$ws.Range("D33").Value = "Igual"
$ws.Range("F33").Value = "premium"
$ws.Range("G33").Value = "https://produto.mercadolivre.com.br/MLB-3470559799-fonte-carregador-automotivo-jfa-60a-bivolt-storm-medidor-cca-_JM#position%3D33%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C34").Value = "FONTE 60A LITE"
$ws.Range("D34").Value = "Igual"
$ws.Range("F34").Value = "premium"
$ws.Range("G34").Value = "https://produto.mercadolivre.com.br/MLB-3350332825-fonte-automotiva-jfa-storm-lite-60a-bivolt-carregador-_JM#position%3D34%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"

$ws.Range("C35").Value = "FONTE 200A LITE"
$ws.Range("D35").Value = "Igual"
$ws.Range("F35").Value = "classico"
$ws.Range("G35").Value = "https://produto.mercadolivre.com.br/MLB-3482419471-fonte-automotiva-jfa-storm-lite-200a-bivolt-carregador-_JM#position%3D35%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D4326e1d9-9a90-4964-80b2-04c03a6bbbdd"
